# Decision tree report: add "Criterion" column, fix a data entry bug in
# rows 11-12, and append new benchmark rows (13-23) for the "entropy"
# criterion run.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Fix a copy/paste data-entry bug in rows 11 and 12: the Training
#    Time / Predition Time / Accuracy values had been shifted one
#    column to the left.
# ---------------------------------------------------------------------
$ws.Range("C11").Value = 265
$ws.Range("D11").Value = 0.23
$ws.Range("E11").Value = 77.52

$ws.Range("C12").Value = 266
$ws.Range("D12").Value = 0.23
$ws.Range("E12").Value = 77.47

# E12 had lost the "Good"/centered/bold styling that the rest of the
# accuracy column in this block (E9:E11) uses - restore it.
$ws.Range("E9").Copy() | Out-Null
$ws.Range("E12").PasteSpecial(-4122) | Out-Null
$ws.Range("E12").Value = 77.47

# ---------------------------------------------------------------------
# 2) Add the new "Criterion" column (G) with a header cell matching the
#    look of the other headers (bold, header-note fill, thin side
#    borders only).
# ---------------------------------------------------------------------
$ws.Range("G1").Value = "Criterion"
$ws.Range("G1").Font.Bold = $true
$ws.Range("G1").Interior.Color = $ws.Range("A1").Interior.Color()
$ws.Range("G1").Borders.Item($xlEdgeLeft).LineStyle = 1
$ws.Range("G1").Borders.Item($xlEdgeLeft).Color = 11711154
$ws.Range("G1").Borders.Item($xlEdgeRight).LineStyle = 1
$ws.Range("G1").Borders.Item($xlEdgeRight).Color = 11711154

# Fill G2:G12 with "gini" (existing rows used the gini criterion),
# style matching the "20% - Accent5" theme used for this column.
for ($r = 2; $r -le 12; $r++) {
    $cell = $ws.Range("G$r")
    $cell.Value = "gini"
    $cell.Style = "20% - Accent5"
}

# ---------------------------------------------------------------------
# 3) Append the new "entropy" benchmark rows (13-23).
# ---------------------------------------------------------------------
$entropyRows = @(
    @{ Row = 13; B = 2;    E = 71.599999999999994;  F = "90 by 100" },
    @{ Row = 14; B = 50;   E = 73.05;                F = "90 by 100" },
    @{ Row = 15; B = 100;  E = 73.430000000000007;  F = "90 by 100" },
    @{ Row = 16; B = 1000; E = 66.03;                F = "90 by 100" },
    @{ Row = 17; B = 500;  E = 73.5;                 F = "90 by 100" },
    @{ Row = 18; B = 250;  E = 73.5;                 F = "90 by 100" },
    @{ Row = 19; B = 300;  E = 72.92;                F = "90 by 100" },
    @{ Row = 20; B = 250;  E = $null;                F = "90 by 50" },
    @{ Row = 21; B = 250;  E = $null;                F = "90 by 10" },
    @{ Row = 22; B = 250;  E = $null;                F = "90 by 5" },
    @{ Row = 23; B = 250;  E = $null;                F = "90 by 5" }
)

foreach ($item in $entropyRows) {
    $r = $item.Row
    $ws.Range("B$r").Value = $item.B

    $eCell = $ws.Range("E$r")
    if ($null -ne $item.E) {
        $eCell.Value = $item.E
    } else {
        $eCell.Value = ""
    }
    $eCell.HorizontalAlignment = -4108

    $ws.Range("F$r").Value = $item.F

    $gCell = $ws.Range("G$r")
    $gCell.Value = "entropy"
    $gCell.Interior.Color = $ws.Range("G2").Interior.Color()
}

# ---------------------------------------------------------------------
# 4) Re-center the existing Accuracy column cells that lost their
#    explicit style index when new styles were introduced (purely
#    cosmetic - keeps E2:E8 centered like before).
# ---------------------------------------------------------------------
foreach ($r in 2,3,4,5,7) {
    $ws.Range("E$r").HorizontalAlignment = -4108
}

# ---------------------------------------------------------------------
# 5) Update the sheet selection to match where the user ended up after
#    entering the new data (cell B23).
# ---------------------------------------------------------------------
$ws.Range("B23").Select() | Out-Null

Write-Host "Decision tree report updated with Criterion column and entropy rows."
